$d = $word.ActiveDocument

# 1. Title paragraph ("IT4472 Programming for AI") loses its default
#    "space after" so it sits flush against the line below it.
$p1 = $d.Paragraphs.Item(1)
$p1.Format.SpaceAfter = 0

# 2. Remove the empty paragraph that used to sit between
#    "Declaration of Use of Artificial Intelligence" and the
#    "Declaration:" heading.
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Delete()

# 3. Append a new "Declaration:" heading plus a paragraph describing the
#    Stack Overflow usage at the very end of the document.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Declaration:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>We also referred to developer discussions and solutions on Stack Overflow</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> [</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>https://stackoverflow.com/questions</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>]</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> while debugging parts of our code and understanding how to apply certain functions. We evaluated multiple threads and selectively applied techniques that were relevant to our implementation. All adaptations were made with full understanding, and we integrated them into our original workflow.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$end = $d.Content.End - 1
$insertionRange = $d.Range($end, $end)
[void]$insertionRange.InsertXML($xml)
